$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update sample data: keep only the header row + a single data row ---
# Remove rows 3 and 4 (the old "Workshop" and "Coaching" sample entries).
$ws.Rows("3:4").Delete()

# Replace the remaining data row (row 2) with the new single sample entry.
$ws.Range("A2").Value = "Designing Effective Moodle Courses for Higher Education"

$description = @"
This professional development course supports higher education teachers in designing, structuring, and implementing pedagogically sound Moodle courses. Participants will learn how to translate didactic concepts into functional Moodle environments that support student engagement, self-regulated learning, and assessment.
The course combines instructional design principles with hands-on practice in Moodle. Participants will explore core Moodle functionalities (e.g., activities, resources, assessments, feedback, and analytics) and learn how to align them with learning objectives, constructive alignment, and evidence-based teaching strategies.
By the end of the course, participants will have developed a prototype Moodle course or a redesigned course unit that is ready for implementation in their own teaching context.
"@
$ws.Range("B2").Value = $description.TrimEnd("`r", "`n")

# --- Style the header row (A1:B1): bold font, thin box border, centered/top aligned ---
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous (thin)
$header.Borders.Weight = 2            # xlThin
